# Apply changes from the commit "push 2 desktop (5/6)"
#
# Sheet "Ship" (sheet1.xml): flip a batch of Truth/Sensed flags between
# 0 and 1 across several row-bands (columns B,C,D,E,F,G correspond to
# Ship Truth, System 1 Truth, System 2 Truth, Ship Sensed, System 1
# Sensed, System 2 Sensed).
#
# Sheet "System 1" (sheet2.xml): append 10 new data rows (rows 4-13)
# with a trailing check formula in column N.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Ship"
# ---------------------------------------------------------------
$ship = $wb.Worksheets.Item("Ship")

$ship.Range("E7:F11").Value  = 0
$ship.Range("B14:C18").Value = 1
$ship.Range("B22:C22").Value = 0
$ship.Range("E26:F31").Value = 1

$ship.Range("B34").Value = 0
$ship.Range("C34").Value = 0
$ship.Range("E34").Value = 0
$ship.Range("F34").Value = 0

$ship.Range("E35:F41").Value = 0
$ship.Range("B42:C51").Value = 1

$ship.Range("B53").Value = 1
$ship.Range("C53").Value = 1
$ship.Range("E53").Value = 1
$ship.Range("F53").Value = 1

$ship.Range("E54:F61").Value = 1
$ship.Range("B63:C71").Value = 1
$ship.Range("B80:C80").Value = 1

$ship.Range("B86:B89").Value = 0
$ship.Range("D86:D89").Value = 0
$ship.Range("B90:D91").Value = 0

$ship.Range("B94:C97").Value = 1

$ship.Range("C98:C101").Value = 1
$ship.Range("D98:D101").Value = 0
$ship.Range("E98:E101").Value = 0
$ship.Range("G98:G101").Value = 0

# ---------------------------------------------------------------
# Sheet "System 1" — append rows 4..13
# ---------------------------------------------------------------
$sys1 = $wb.Worksheets.Item("System 1")

$rows = @(
    @(2,  1,1,1,1,1,1,1,1,1,1,1,1),
    @(3,  1,1,1,1,1,1,1,1,1,1,1,1),
    @(4,  1,1,1,1,1,1,1,1,1,1,1,1),
    @(5,  0,1,1,1,1,1,0,1,1,1,1,1),
    @(6,  0,1,1,1,1,1,0,1,1,1,1,1),
    @(7,  0,1,1,1,1,1,0,1,1,1,1,1),
    @(8,  0,1,1,1,1,1,0,1,1,1,1,1),
    @(9,  0,1,1,1,1,1,0,1,1,1,1,1),
    @(10, 1,1,1,1,1,1,1,1,1,1,1,1),
    @(11, 1,1,1,1,1,1,1,1,1,1,1,1)
)

$r = 4
foreach ($row in $rows) {
    $sys1.Cells.Item($r, 1).Value = $row[0]          # A: Time Step
    for ($col = 2; $col -le 13; $col++) {
        $sys1.Cells.Item($r, $col).Value = $row[$col - 1]   # B..M: component truth/sensed flags
    }
    $sys1.Cells.Item($r, 14).Formula = "=IF(B$r = H$r, 1, 0)"  # N: check column
    $r++
}
